# correção das notas do fórum para matc65 em 2021.2
# For every row where column J (nota_view) currently equals 4,
# zero out the forum-view flags (B:H), total_views (I) and nota_view (J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $notaView = $ws.Cells.Item($r, 10).Value()
    if ($notaView -eq 4) {
        for ($c = 2; $c -le 10; $c++) {
            $ws.Cells.Item($r, $c).Value = 0
        }
    }
}
